$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (D) and "Volumen" (M) values between rows 3 and 4.
$d3 = $ws.Range("D3").Value2
$d4 = $ws.Range("D4").Value2
$m3 = $ws.Range("M3").Value2
$m4 = $ws.Range("M4").Value2

$ws.Range("D3").Value2 = $d4
$ws.Range("D4").Value2 = $d3

$ws.Range("M3").Value2 = $m4
$ws.Range("M4").Value2 = $m3
